$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-454) holds a "Förändrad" (changed) date serial that was
# bumped from 45175 (2023-09-06) to 45177 (2023-09-08) for every data row.
$oldValue = 45175
$newValue = 45177
$lastRow = 454

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
